# "aggiornamento a 9/09 compreso" - append daily rows for 2021-09-02 .. 2021-09-09
# (rows 367..374, continuing the existing A:D data table) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily records: row, date-serial (col A), nuovi pos. (col B),
# somma mobile 7gg. (col C), somma mobile 7gg. per 100mila abitanti (col D)
$data = @(
    @(367, 44441, 0, 2, 96.15384615384616),
    @(368, 44442, 0, 1, 48.07692307692308),
    @(369, 44443, 0, 1, 48.07692307692308),
    @(370, 44444, 0, 0, 0),
    @(371, 44445, 0, 0, 0),
    @(372, 44446, 0, 0, 0),
    @(373, 44447, 0, 0, 0),
    @(374, 44448, 0, 0, 0)
)

# Column A carries the bold/centered/bordered date style (same as every
# preceding row in the table) - copy it down from the last existing row
# before writing the new values so formatting stays consistent.
$lastRow = 366
$styleSrc = $ws.Cells.Item($lastRow, 1)

foreach ($row in $data) {
    $r = $row[0]

    $styleSrc.Copy($ws.Cells.Item($r, 1))
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

Write-Output "Appended rows 367-374 to $($ws.Name)"
